$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh aggregate stats now that trade #8 has closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.09   # Current Capital
$summary.Range("B4").Value = 0.09      # Total P&L $
$summary.Range("B5").Value = 0.23      # Total P&L %
$summary.Range("B6").Value = 8         # Total Trades
$summary.Range("B7").Value = 4         # Winning Trades
$summary.Range("B9").Value = 50        # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 4) mirrors the same refresh.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.09     # Capital
$status.Range("D4").Value = 8          # Trades
$status.Range("E4").Value = 0.09       # P&L $
$status.Range("F4").Value = 0.09       # P&L %
$status.Range("G4").Value = 50         # Win Rate %

# ---------------------------------------------------------------------------
# Helper: append the newly-closed trade #8 row to a trade log sheet.
# ---------------------------------------------------------------------------
function Add-Trade8Row($sheet) {
    $sheet.Range("A9").Value = 8

    # Date/Time columns hold literal text in this workbook (not Excel date
    # serials) - force text formatting before assignment, then restore the
    # default style so no stray number-format style sticks to the cell.
    $sheet.Range("B9").NumberFormat = "@"
    $sheet.Range("B9").Value = "2026-02-17"
    $sheet.Range("B9").Style = "Normal"

    $sheet.Range("C9").Value = "08:08:12"
    $sheet.Range("D9").Value = "MarketMaking"
    $sheet.Range("E9").Value = "UP"
    $sheet.Range("F9").Value = 0.26
    $sheet.Range("G9").Value = 0.43
    $sheet.Range("H9").Value = "CLOSED"
    $sheet.Range("I9").Value = 65.38460000000001
    $sheet.Range("J9").Value = 0.17
    $sheet.Range("K9").Value = 100.09
    $sheet.Range("L9").Value = 0
    $sheet.Range("M9").Value = 0
    $sheet.Range("N9").Value = 0.6
    $sheet.Range("O9").Value = "Normal spread capture: 19600 bps"
    $sheet.Range("P9").Value = "early_exit"
    $sheet.Range("Q9").Value = 0.13
}

Add-Trade8Row $wb.Worksheets.Item("All Trades")
Add-Trade8Row $wb.Worksheets.Item("MarketMaking")
